# Implements: "Implemented logic to skip step definitions based on
# parameter given in testdata"
#
# Target sheet: "NitroXHome" (4th tab / sheet4.xml). Adds two new columns to
# the right of the existing A:L test-data table:
#   M -> "SkipAtStepNum"  (header)
#   N -> "SIT"             (header), data rows default to boolean FALSE
# and centers (+ keeps wrap/border) the whole header row, matching the
# existing bordered/wrap-text look already used by columns A-C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # Amazon, RestAPITestData, NitroXLogin, NitroXHome

# ---------------------------------------------------------------------
# 1) Give the previously-unstyled header cells (D1:N1) and the new M/N
#    data columns the same bordered/wrap-text base format already used
#    by the existing header cells (A1:C1), by copying C1's format over.
#    Doing this *before* centering lets Excel collapse everything down to
#    the same two new cell styles instead of a combinatorial explosion.
# ---------------------------------------------------------------------
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:N13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Center-align the header row and the new M/N columns (wrap + border
#    are already in place from the paste above).
# ---------------------------------------------------------------------
$ws.Range("A1:N1").HorizontalAlignment = -4108        # xlCenter
$ws.Range("M2:M13").HorizontalAlignment = -4108
$ws.Range("N2:N13").HorizontalAlignment = -4108

# M10 keeps a Text number format (as in the source workbook).
$ws.Range("M10").NumberFormat = "@"

# ---------------------------------------------------------------------
# 3) New column widths.
# ---------------------------------------------------------------------
$ws.Columns("M").ColumnWidth = 17.4

# ---------------------------------------------------------------------
# 4) Header text. "SIT" is written before "SkipAtStepNum" so the shared
#    string table picks up the same ordering as the target workbook.
# ---------------------------------------------------------------------
$ws.Range("N1").Value = "SIT"
$ws.Range("M1").Value = "SkipAtStepNum"

# ---------------------------------------------------------------------
# 5) Data rows: SIT column defaults to FALSE for every existing row.
# ---------------------------------------------------------------------
$ws.Range("N2:N13").Value = $false

# ---------------------------------------------------------------------
# 6) Selection state, matching the saved workbook view.
# ---------------------------------------------------------------------
$ws.Range("F10:L11").Select() | Out-Null

# ---------------------------------------------------------------------
# 7) Misc window-state cosmetics (best-effort; not all UI chrome is
#    reachable through automation in a headless runtime).
# ---------------------------------------------------------------------
try { $wb.Windows.Item(1).Height = 10300 } catch {}
